$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '61.187.06'
$c.ClearFormats()
$c = $ws.Range('E2')
$c.NumberFormat = '@'
$c.Value = '  -4.85%  '
$c.ClearFormats()
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '3.315.65'
$c.ClearFormats()
$c = $ws.Range('E3')
$c.NumberFormat = '@'
$c.Value = '  -5.12%  '
$c.ClearFormats()
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.ClearFormats()
$c = $ws.Range('E4')
$c.NumberFormat = '@'
$c.Value = '  +0.03%  '
$c.ClearFormats()
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '569.89'
$c.ClearFormats()
$c = $ws.Range('E5')
$c.NumberFormat = '@'
$c.Value = '  -3.18%  '
$c.ClearFormats()
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '126.81'
$c.ClearFormats()
$c = $ws.Range('E6')
$c.NumberFormat = '@'
$c.Value = '  -5.33%  '
$c.ClearFormats()
$c = $ws.Range('E7')
$c.NumberFormat = '@'
$c.Value = '  -0.05%  '
$c.ClearFormats()
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '3.316.89'
$c.ClearFormats()
$c = $ws.Range('E8')
$c.NumberFormat = '@'
$c.Value = '  -5.04%  '
$c.ClearFormats()
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.477'
$c.ClearFormats()
$c = $ws.Range('E9')
$c.NumberFormat = '@'
$c.Value = '  -2.18%  '
$c.ClearFormats()
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '7.21'
$c.ClearFormats()
$c = $ws.Range('E10')
$c.NumberFormat = '@'
$c.Value = '  -5.12%  '
$c.ClearFormats()
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.118'
$c.ClearFormats()
$c = $ws.Range('E11')
$c.NumberFormat = '@'
$c.Value = '  -5.51%  '
$c.ClearFormats()
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.374'
$c.ClearFormats()
$c = $ws.Range('E12')
$c.NumberFormat = '@'
$c.Value = '  -4.07%  '
$c.ClearFormats()
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '3.878.23'
$c.ClearFormats()
$c = $ws.Range('E14')
$c.NumberFormat = '@'
$c.Value = '  -1.46%  '
$c.ClearFormats()
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '3.313.90'
$c.ClearFormats()
$c = $ws.Range('E15')
$c.NumberFormat = '@'
$c.Value = '  -5.21%  '
$c.ClearFormats()
$c = $ws.Range('E16')
$c.NumberFormat = '@'
$c.Value = '  -6.87%  '
$c.ClearFormats()
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '24.72'
$c.ClearFormats()
$c = $ws.Range('E17')
$c.NumberFormat = '@'
$c.Value = '  -2.72%  '
$c.ClearFormats()
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '61.259.65'
$c.ClearFormats()
$c = $ws.Range('E18')
$c.NumberFormat = '@'
$c.Value = '  -4.75%  '
$c.ClearFormats()
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '5.60'
$c.ClearFormats()
$c = $ws.Range('E19')
$c.NumberFormat = '@'
$c.Value = '  -3.04%  '
$c.ClearFormats()
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '9.08'
$c.ClearFormats()
$c = $ws.Range('E20')
$c.NumberFormat = '@'
$c.Value = '  -9.34%  '
$c.ClearFormats()
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '13.18'
$c.ClearFormats()
$c = $ws.Range('E21')
$c.NumberFormat = '@'
$c.Value = '  -2.52%  '
$c.ClearFormats()
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '352.97'
$c.ClearFormats()
$c = $ws.Range('E22')
$c.NumberFormat = '@'
$c.Value = '  -9.11%  '
$c.ClearFormats()
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '0.554'
$c.ClearFormats()
$c = $ws.Range('E23')
$c.NumberFormat = '@'
$c.Value = '  -4.60%  '
$c.ClearFormats()
$c = $ws.Range('E24')
$c.NumberFormat = '@'
$c.Value = '  -0.02%  '
$c.ClearFormats()
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '3.444.97'
$c.ClearFormats()
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '70.19'
$c.ClearFormats()
$c = $ws.Range('E26')
$c.NumberFormat = '@'
$c.Value = '  -5.40%  '
$c.ClearFormats()
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '0.0000107'
$c.ClearFormats()
$c = $ws.Range('E27')
$c.NumberFormat = '@'
$c.Value = '  -7.01%  '
$c.ClearFormats()
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.ClearFormats()
$c = $ws.Range('E28')
$c.NumberFormat = '@'
$c.Value = '  +0.11%  '
$c.ClearFormats()
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '7.14'
$c.ClearFormats()
$c = $ws.Range('E29')
$c.NumberFormat = '@'
$c.Value = '  -3.31%  '
$c.ClearFormats()
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '1.47'
$c.ClearFormats()
$c = $ws.Range('E30')
$c.NumberFormat = '@'
$c.Value = '  -2.85%  '
$c.ClearFormats()
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '7.89'
$c.ClearFormats()
$c = $ws.Range('E31')
$c.NumberFormat = '@'
$c.Value = '  -3.52%  '
$c.ClearFormats()
$c = $ws.Range('E32')
$c.NumberFormat = '@'
$c.Value = '  -6.61%  '
$c.ClearFormats()
$c = $ws.Range('E33')
$c.NumberFormat = '@'
$c.Value = '  -0.08%  '
$c.ClearFormats()
$c = $ws.Range('E34')
$c.NumberFormat = '@'
$c.Value = '  -5.04%  '
$c.ClearFormats()
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '3.340.44'
$c.ClearFormats()
$c = $ws.Range('E35')
$c.NumberFormat = '@'
$c.Value = '  -5.20%  '
$c.ClearFormats()
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '22.34'
$c.ClearFormats()
$c = $ws.Range('E36')
$c.NumberFormat = '@'
$c.Value = '  -4.27%  '
$c.ClearFormats()
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '5.43'
$c.ClearFormats()
$c = $ws.Range('E37')
$c.NumberFormat = '@'
$c.Value = '  +1.77%  '
$c.ClearFormats()
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '6.73'
$c.ClearFormats()
$c = $ws.Range('E38')
$c.NumberFormat = '@'
$c.Value = '  -2.77%  '
$c.ClearFormats()
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '162.76'
$c.ClearFormats()
$c = $ws.Range('E39')
$c.NumberFormat = '@'
$c.Value = '  -1.75%  '
$c.ClearFormats()
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '1.49'
$c.ClearFormats()
$c = $ws.Range('E40')
$c.NumberFormat = '@'
$c.Value = '  -3.57%  '
$c.ClearFormats()
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.0754'
$c.ClearFormats()
$c = $ws.Range('E41')
$c.NumberFormat = '@'
$c.Value = '  -4.19%  '
$c.ClearFormats()
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.ClearFormats()
$c = $ws.Range('E42')
$c.NumberFormat = '@'
$c.Value = '  -0.08%  '
$c.ClearFormats()
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '40.94'
$c.ClearFormats()
$c = $ws.Range('E43')
$c.NumberFormat = '@'
$c.Value = '  -2.59%  '
$c.ClearFormats()
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '0.748'
$c.ClearFormats()
$c = $ws.Range('E44')
$c.NumberFormat = '@'
$c.Value = '  -7.48%  '
$c.ClearFormats()
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '4.24'
$c.ClearFormats()
$c = $ws.Range('E45')
$c.NumberFormat = '@'
$c.Value = '  -4.12%  '
$c.ClearFormats()
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '1.12'
$c.ClearFormats()
$c = $ws.Range('E46')
$c.NumberFormat = '@'
$c.Value = '  -5.51%  '
$c.ClearFormats()
$c = $ws.Range('E47')
$c.NumberFormat = '@'
$c.Value = '  -6.49%  '
$c.ClearFormats()
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '22.47'
$c.ClearFormats()
$c = $ws.Range('E48')
$c.NumberFormat = '@'
$c.Value = '  -8.31%  '
$c.ClearFormats()
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '6.63'
$c.ClearFormats()
$c = $ws.Range('E49')
$c.NumberFormat = '@'
$c.Value = '  -2.66%  '
$c.ClearFormats()
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.851'
$c.ClearFormats()
$c = $ws.Range('E50')
$c.NumberFormat = '@'
$c.Value = '  -7.40%  '
$c.ClearFormats()
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '2.195.81'
$c.ClearFormats()
$c = $ws.Range('E51')
$c.NumberFormat = '@'
$c.Value = '  -8.44%  '
$c.ClearFormats()
